# Add the new "Sebastiano Zoller" team as row 61 (the sheet currently ends at row 60).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A61").Value = "Sebastiano Zoller"
$ws.Range("B61").Value = "Lorenzo Canali | CGB Gamberoni"
$ws.Range("C61").Value = "Sebastiano Zoller | CGB Gamberoni"
$ws.Range("D61").Value = "Andrea  Roveda  | Pinguini Trentini"
$ws.Range("E61").Value = "Michele Merighi | Clitoriders"
$ws.Range("F61").Value = "Alessio Debiasi | Mai una gioia"
